$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new errata row (Chap 24, p.852, 24.3.3, "spread" -> "spreads") is being
# inserted just above the existing "Chap 24 / p.871 / 24.6" row. Rather than
# shifting every subsequent row down (which would renumber row 81 etc.), the
# existing row 77 contents are pushed down one row into (until-now empty)
# row 78, and the new entry is written into row 77 in its place.

# Step 1: carry the old row 77 (B:G) data down into row 78 (column A is left
# blank there, matching the rest of the sheet's "continuation row" layout).
$ws.Range("B78").Value = 871
$ws.Range("C78").Value = 24.6
$ws.Range("D78").Value = "last"
$ws.Range("E78").Value = 3
$ws.Range("F78").Value = "it it"
$ws.Range("G78").Value = "if it"

# Step 2: write the new errata entry into row 77 (column A, "Chap 24",
# already holds the right value and is left untouched).
$ws.Range("B77").Value = 852
$ws.Range("C77").Value = "24.3.3"
$ws.Range("D77").Value = 2
$ws.Range("E77").Value = 1
$ws.Range("F77").Value = "spread"
$ws.Range("G77").Value = "spreads"

# Update the view: scrolled down a bit further and a new active cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 54
$win.ScrollColumn = 1
$ws.Range("H77").Select()
